# Add a "packages" worksheet (second-page sidebar) listing the author's
# R packages, matching the style already used by the other reference
# sheets (languages, software, ...).

$wb = $excel.ActiveWorkbook

# --- Cosmetic tidy-up of the existing "languages" sheet (done first so the
#     new "packages" sheet ends up the active / selected tab at the end). ---
$langWs = $wb.Worksheets.Item("languages")
$langWs.Columns.Item(1).ColumnWidth = 14.498697916666666
$langWs.Columns.Item(2).ColumnWidth = 23.166666666666668
$langWs.Activate()
$langWs.Rows("1:6").Select()

# --- Add the new "packages" sheet after the last existing sheet. ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "packages"

# Column widths matching the other metadata sheets.
$ws.Columns.Item(1).ColumnWidth = 14.498697916666666
$ws.Columns.Item(2).ColumnWidth = 17.333333333333332

# --- Header row (documentation row, like the other sheets). ---
$ws.Range("A1").Value = "Name of package"
$ws.Range("C1").Value = "Link"
$ws.Range("B1").Value = "Description of the package"

# --- Meta/label row. ---
$ws.Range("B2").Value = "description   "

# --- Package rows. ---
$ws.Range("A3").Value = "dangeo"
$ws.Range("C3").Value = "https://github.com/javiereliomedina/dangeo.git"

$ws.Range("A4").Value = "ggpyramid"
$ws.Range("C4").Value = "https://github.com/javiereliomedina/ggpyramid.git"

$ws.Range("A5").Value = "ggrugby"
$ws.Range("C5").Value = "https://github.com/javiereliomedina/ggrugby.git"

$ws.Range("A6").Value = "safecastR"
$ws.Range("C6").Value = "https://github.com/javiereliomedina/safecastR.git"

$ws.Range("B3").Value = "R package for accessing the Danish Map Supply Download via the kortforsyningen FTP"
$ws.Range("B4").Value = "R function for plotting population pyramids in {ggplot2}"
$ws.Range("B5").Value = "R functions for plotting rugby events in {ggplot2}"
$ws.Range("B6").Value = "R package for loading data from Safecast API"

$ws.Range("A2").Value = "package"
$ws.Range("C2").Value = "link"

# Bold/style row 4 (dangeo/ggpyramid? no - matches the "in_resume" style row
# used on other sheets) the same way the source sheet uses style index 1 on
# row 4 of the languages sheet.
$ws.Range("A4:B4").Style = $langWs.Range("A4:B4").Style

# --- Turn the last URL into a real hyperlink (adds the Hyperlink style). ---
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/javiereliomedina/safecastR.git")

# --- Final selection / active sheet state. ---
$ws.Range("A8").Select()
$ws.Activate()
